# actualizacion 05 junio de 2016
# Column B holds ISO-8601 datetime-with-offset strings (e.g.
# "2016-05-16T19:02:53+02:00"). Replace each with just the time-of-day
# portion (e.g. "19:02:53"), dropping the date and the UTC offset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2

    if ($val -is [string] -and $val -match "^\d{4}-\d{2}-\d{2}T\d{2}:\d{2}:\d{2}") {
        $timePart = $val.Split("T")[1]
        $timePart = $timePart.Split("+")[0]
        $timePart = $timePart.Split("-")[0]
        $cell.Value = $timePart
    }
}
